$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.123.98"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.814.46"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.80"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4409"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3727"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.66"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07687"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.115"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.89"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.297"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.500"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "1.827.15"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.04"
$ws.Range("E17").Value = "  +11.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001081"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06469"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.51"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.282"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5384"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "28.184.10"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085"
$ws.Range("E26").Value = "  -13.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.59"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "2.030.50"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.328"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.45"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.189"
$ws.Range("E32").Value = "  -10.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.831"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09245"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.665"
$ws.Range("E35").Value = "  -7.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.95"
$ws.Range("E36").Value = "  +5.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02328"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2172"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.147"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6540"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06157"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.195"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.098"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.391"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.81"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6049"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.764"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.47"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06978"
$ws.Range("E51").Value = "  -0.07%  "
